$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Row -> (PeriodToExpire, LastUpdateDate)
$updates = @{
    3  = @(164, "04-Nov-2025")
    4  = @(310, "04-Nov-2025")
    5  = @(355, "04-Nov-2025")
    6  = @(355, "04-Nov-2025")
    7  = @(161, "04-Nov-2025")
    8  = @(304, "04-Nov-2025")
    9  = @(-34, "04-Nov-2025")
    10 = @(-23, "04-Nov-2025")
    11 = @(630, "04-Nov-2025")
    12 = @(630, "04-Nov-2025")
    13 = @(630, "04-Nov-2025")
    14 = @(630, "04-Nov-2025")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 8).Value = $vals[0]
    $dateCell = $ws.Cells.Item($row, 9)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $vals[1]
}
